$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RunManager")
$ws2 = $wb.Worksheets.Item("TestData")

$ws1.Activate()
$ws1.Range("D4").Select()

$ws2.Range("D1").Value = "valueforsearch"
$ws2.Range("D2").Value = "automation"
$ws2.Range("D3").Value = "testing"

$ws2.Range("A4").Value = "test1"
$ws2.Range("B4").Value = "abcd"
$ws2.Range("C4").Value = "1234"
$ws2.Range("D4").Value = "test1"

$ws2.Range("A5").Value = "test1"
$ws2.Range("B5").Value = "efgh"
$ws2.Range("C5").Value = "4567"
$ws2.Range("D5").Value = "test11"

$ws2.Range("A6").Value = "test2"
$ws2.Range("B6").Value = "jbvb"
$ws2.Range("C6").Value = "987"
$ws2.Range("D6").Value = "selenium"

$ws2.Range("A7").Value = "test3"
$ws2.Range("B7").Value = "hgb1"
$ws2.Range("C7").Value = "jnh"
$ws2.Range("D7").Value = "appium"

$ws2.Activate()
$ws2.Range("A4").Select()
